$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
